$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 18:46"

# --- Country list was re-sorted (by total cases); update the country name
#     shown in each affected row (A column) so it matches the new order ---
$names = @{
    48 = "Catar"
    49 = "Eslovenia"
    71 = "Marruecos"
    72 = "Letonia"
    73 = "Eslovaquia"
    74 = "Nueva Zelanda"
    75 = "Kuwait"
    76 = "Uruguay"
    77 = "Principado de Andorra"
    78 = "San Marino"
    79 = "Costa Rica"
    80 = "Republica de Macedonia"
    81 = "Tunez"
    82 = "Jordania"
    83 = "Bosnia y Herzegovina"
}

foreach ($r in $names.Keys) {
    $ws.Cells.Item([int]$r, 1).Value = $names[$r]
}

# --- Update statistics (columns B..H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#     for every row whose numbers changed in this update ---
$rowData = @{
    11 = @(10897, 1020, 131, 10613, 141, 31, 153)
    24 = @(1872, 0, 26, 1802, 136, 0, 44)
    48 = @(537, 11, 41, 496, 6, 0, 0)
    49 = @(528, 48, 10, 513, 14, 1, 5)
    71 = @(225, 55, 7, 212, 1, 1, 6)
    72 = @(221, 24, 1, 220, 0, 0, 0)
    73 = @(216, 12, 7, 209, 2, 0, 0)
    74 = @(205, 0, 22, 183, 0, 0, 0)
    75 = @(195, 4, 43, 152, 6, 0, 0)
    76 = @(189, 0, 0, 189, 3, 0, 0)
    77 = @(188, 24, 1, 186, 6, 0, 1)
    78 = @(187, 0, 4, 162, 12, 0, 21)
    79 = @(177, 0, 2, 173, 4, 0, 2)
    80 = @(177, 29, 1, 173, 1, 1, 3)
    81 = @(173, 59, 2, 166, 11, 1, 5)
    82 = @(172, 18, 1, 171, 0, 0, 0)
    83 = @(168, 0, 2, 163, 1, 0, 3)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item([int]$r, $c + 2).Value = $vals[$c]
    }
}
